$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "The reel features a mix of emotional and romantic content with high engagement, indicating a strong connection with the audience."
$ws.Range("B7").Value = "Romantic and emotional content creator with a focus on love and relationships."
$ws.Range("C7").Value = "Positive and engaging, with a lot of love and emotional reactions."
$ws.Range("D7").Value = "High"
$ws.Range("E7").Value = "both"
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = 18500
$ws.Range("H7").Value = "https://www.instagram.com/reels/DIbGiMcxUvJ/"
$ws.Range("I7").Value = "https://www.instagram.com/_pihu_arya2425/"
$ws.Range("J7").Value = "The reel has high engagement and a romantic theme, which aligns well with Knytt's focus on personal connections. The creator's audience is likely to be interested in a texting & video calling app."
$ws.Range("K7").Value = "both"
$ws.Range("L7").Value = "Hey Pihu! Loved your reel—such a beautiful expression of love! 💖 We’re building Knytt, a cozy space for heartfelt conversations. Would love for you to check it out and maybe share your thoughts with your audience. Let’s connect! 😊"
$ws.Range("M7").Value = "This is so heartfelt! 💖 If you ever want to share more love stories, Knytt is perfect for deep, personal convos. Check it out! 😊"
